# Apply the "New template for CovidToscana.xlsx" edit.
#
# Functional changes reproduced here:
#   1. Workbook-level defined names pointing at the Template summary row.
#   2. Two new helper rows (25/26) on the "Template" sheet that pick out a
#      single province's values (via OFFSET + a free "province index" cell)
#      so the pie chart can show one province instead of a fixed range.
#   3. The "Grafico 4" pie chart's Values series is repointed from the old
#      (now unused) Template!D14:G14 range onto the new I26:L26 helper cells.
#   4. Restore the active-cell selection on the Template sheet to E4 (the
#      state the workbook was left in after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# --- 1. Defined names -------------------------------------------------
$wb.Names.Add("ActiveCases", '=Template!$G$4')
$wb.Names.Add("DailyCases", '=Template!$B$4')
$wb.Names.Add("Deaths", '=Template!$E$4')
$wb.Names.Add("GlobalCases", '=Template!$C$4')
$wb.Names.Add("IntensiveCases", '=Template!$F$4')
$wb.Names.Add("Recovered", '=Template!$D$4')

# --- 2. New helper rows 25 & 26 on the Template sheet ------------------
# Entry order matters for shared-string interning: "Num prov:" / "Valori:"
# are entered first (rows 25/26 col H), then the "${numProvince}" token.
$ws.Range("H25").Value = "Num prov:"
$ws.Range("H26").Value = "Valori:"
$ws.Range("I25").Value = '${numProvince}'

$ws.Range("I26").Formula = '=OFFSET(D4, I25, 0)'
$ws.Range("J26").Formula = '=OFFSET(E4, I25, 0)'
$ws.Range("K26").Formula = '=OFFSET(F4, I25, 0)'
$ws.Range("L26").Formula = '=OFFSET(G4, I25, 0)'

# --- 3. Repoint the pie chart ("Grafico 4" / chart2.xml) onto the new
#        helper cells instead of the old Template!D14:G14 range.
$co = $ws.ChartObjects(2)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.Values = '(Template!$I$26,Template!$J$26,Template!$K$26,Template!$L$26)'

# --- 4. Restore selection state on the Template sheet ------------------
$ws.Activate() | Out-Null
$ws.Range("E4").Select() | Out-Null
